$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.753.47'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.077.10'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.51'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.41'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0783'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.384.57'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.75'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.772'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.078.05'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.717.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.08'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.23'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.37'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.40'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.138'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.03'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.46'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.39'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.69'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.34%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.66'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.29%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.41'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.32'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('B40').Value = 'Cronos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0975'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.76'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.23%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.449.84'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.16'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.40'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.58%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.43'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.268.91'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.39%  '
